$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = [double]"22.36000000000006"
$ws.Range("H2").Value = [double]"4.825806421138168e-12"
$ws.Range("I2").Value = [double]"4.825806421138168e-12"
$ws.Range("L2").Value = [double]"54.47784873262754"
$ws.Range("M2").Value = "[42.76133595832469, 66.19436150693039]"
$ws.Range("N2").Value = [double]"3.897770994854e-12"
$ws.Range("O2").Value = [double]"3.897770994854e-12"
$ws.Range("P2").Value = [double]"1.540921321580579"
$ws.Range("Q2").Value = "[1.3019212798660407, 1.7799213632951174]"
$ws.Range("T2").Value = [double]"48.24900798820318"
$ws.Range("U2").Value = "[40.65043580600482, 55.84758017040154]"
$ws.Range("X2").Value = [double]"16.87631631631636"
$ws.Range("Y2").Value = [double]"16.02578578578583"
$ws.Range("Z2").Value = [double]"17.72684684684689"
$ws.Range("F3").Value = [double]"22.36000000000006"
$ws.Range("H3").Value = [double]"1.114097702981098e-11"
$ws.Range("I3").Value = [double]"1.114097702981098e-11"
$ws.Range("L3").Value = [double]"60.59590422484739"
$ws.Range("M3").Value = "[44.352110301961, 76.83969814773377]"
$ws.Range("N3").Value = [double]"1.776617741811037e-09"
$ws.Range("O3").Value = [double]"1.776617741811037e-09"
$ws.Range("P3").Value = [double]"2.018921405009658"
$ws.Range("Q3").Value = "[1.729605565039427, 2.308237244979888]"
$ws.Range("T3").Value = [double]"53.34285758876469"
$ws.Range("U3").Value = "[44.484244528968304, 62.20147064856107]"
$ws.Range("V3").Value = [double]"8.881784197001252e-16"
$ws.Range("W3").Value = [double]"8.881784197001252e-16"
$ws.Range("X3").Value = [double]"15.17525525525529"
$ws.Range("Y3").Value = [double]"14.1456656656657"
$ws.Range("Z3").Value = [double]"16.20484484484488"
$ws.Range("F4").Value = [double]"22.36000000000006"
$ws.Range("H4").Value = [double]"2.875477633779155e-14"
$ws.Range("I4").Value = [double]"2.875477633779155e-14"
$ws.Range("L4").Value = [double]"62.92851256737306"
$ws.Range("M4").Value = "[47.653596279744995, 78.20342885500114]"
$ws.Range("N4").Value = [double]"1.276780903225472e-10"
$ws.Range("O4").Value = [double]"1.276780903225472e-10"
$ws.Range("P4").Value = [double]"2.421447791055197"
$ws.Range("Q4").Value = "[2.1824477493406578, 2.6604478327697354]"
$ws.Range("T4").Value = [double]"54.84653794818048"
$ws.Range("U4").Value = "[47.022980867500955, 62.67009502886001]"
$ws.Range("X4").Value = [double]"13.74278278278281"
$ws.Range("Y4").Value = [double]"12.89225225225228"
$ws.Range("Z4").Value = [double]"14.59331331331335"
$ws.Range("F5").Value = [double]"22.36000000000006"
$ws.Range("H5").Value = [double]"4.611311332780588e-12"
$ws.Range("I5").Value = [double]"4.611311332780588e-12"
$ws.Range("L5").Value = [double]"59.08512008648854"
$ws.Range("M5").Value = "[45.05539168642629, 73.11484848655078]"
$ws.Range("N5").Value = [double]"6.921663242565046e-11"
$ws.Range("O5").Value = [double]"6.921663242565046e-11"
$ws.Range("P5").Value = [double]"2.974921571867812"
$ws.Range("Q5").Value = "[2.698184681461504, 3.25165846227412]"
$ws.Range("T5").Value = [double]"58.48423234371752"
$ws.Range("U5").Value = "[49.73386096484552, 67.23460372258953]"
$ws.Range("X5").Value = [double]"11.77313313313316"
$ws.Range("Y5").Value = [double]"10.78830830830834"
$ws.Range("Z5").Value = [double]"12.75795795795799"
$ws.Range("F6").Value = [double]"24.78000000000043"
$ws.Range("H6").Value = [double]"2.238209617644316e-13"
$ws.Range("I6").Value = [double]"2.238209617644316e-13"
$ws.Range("L6").Value = [double]"60.25739061172973"
$ws.Range("M6").Value = "[48.20220157537892, 72.31257964808053]"
$ws.Range("N6").Value = [double]"4.212186155427844e-13"
$ws.Range("O6").Value = [double]"4.212186155427844e-13"
$ws.Range("P6").Value = [double]"-3.094421592725082"
$ws.Range("Q6").Value = "[-3.308263735311774, -2.8805794501383892]"
$ws.Range("T6").Value = [double]"55.09852233971888"
$ws.Range("U6").Value = "[47.374715867549504, 62.82232881188826]"
$ws.Range("V6").Value = [double]"0"
$ws.Range("W6").Value = [double]"0"
$ws.Range("X6").Value = [double]"12.20396396396418"
$ws.Range("Y6").Value = [double]"11.3606006006008"
$ws.Range("Z6").Value = [double]"13.04732732732756"
$ws.Range("F7").Value = [double]"24.78000000000043"
$ws.Range("H7").Value = [double]"3.93751697913558e-12"
$ws.Range("I7").Value = [double]"3.93751697913558e-12"
$ws.Range("L7").Value = [double]"53.91935573747271"
$ws.Range("M7").Value = "[41.112094081415165, 66.72661739353026]"
$ws.Range("N7").Value = [double]"6.984701705903262e-11"
$ws.Range("O7").Value = [double]"6.984701705903262e-11"
$ws.Range("P7").Value = [double]"2.974921571867812"
$ws.Range("Q7").Value = "[2.7107636310254275, 3.2390795127101963]"
$ws.Range("T7").Value = [double]"55.32556749357759"
$ws.Range("U7").Value = "[47.737109228409324, 62.914025758745865]"
$ws.Range("V7").Value = [double]"0"
$ws.Range("W7").Value = [double]"0"
$ws.Range("X7").Value = [double]"13.04732732732756"
$ws.Range("Y7").Value = [double]"12.00552552552574"
$ws.Range("Z7").Value = [double]"14.08912912912938"
$ws.Range("F8").Value = [double]"24.78000000000043"
$ws.Range("H8").Value = [double]"1.124511594952082e-11"
$ws.Range("I8").Value = [double]"1.124511594952082e-11"
$ws.Range("L8").Value = [double]"62.30722883681198"
$ws.Range("M8").Value = "[45.50660694267165, 79.10785073095231]"
$ws.Range("N8").Value = [double]"2.061335546699183e-09"
$ws.Range("O8").Value = [double]"2.061335546699183e-09"
$ws.Range("Q8").Value = "[2.4969214884387347, 3.0252373701235036]"
$ws.Range("T8").Value = [double]"54.80388398002343"
$ws.Range("U8").Value = "[45.675662821685165, 63.93210513836169]"
$ws.Range("V8").Value = [double]"1.110223024625157e-15"
$ws.Range("W8").Value = [double]"1.110223024625157e-15"
$ws.Range("X8").Value = [double]"13.89069069069093"
$ws.Range("Y8").Value = [double]"12.84888888888912"
$ws.Range("Z8").Value = [double]"14.93249249249275"
$ws.Range("F9").Value = [double]"24.78000000000043"
$ws.Range("H9").Value = [double]"6.560907372943348e-11"
$ws.Range("I9").Value = [double]"6.560907372943348e-11"
$ws.Range("L9").Value = [double]"57.66744986638336"
$ws.Range("M9").Value = "[41.55244740693908, 73.78245232582765]"
$ws.Range("N9").Value = [double]"5.022730453774216e-09"
$ws.Range("O9").Value = [double]"5.022730453774216e-09"
$ws.Range("P9").Value = [double]"2.710763631025428"
$ws.Range("Q9").Value = "[2.408868841491273, 3.012658420559582]"
$ws.Range("T9").Value = [double]"54.57640694735886"
$ws.Range("U9").Value = "[45.591795470196786, 63.56101842452094]"
$ws.Range("V9").Value = [double]"6.661338147750939e-16"
$ws.Range("W9").Value = [double]"6.661338147750939e-16"
$ws.Range("X9").Value = [double]"14.08912912912938"
$ws.Range("Y9").Value = [double]"12.89849849849872"
$ws.Range("Z9").Value = [double]"15.27975975976003"
$ws.Range("F10").Value = [double]"24.78000000000043"
$ws.Range("H10").Value = [double]"8.286704655802168e-13"
$ws.Range("I10").Value = [double]"8.286704655802168e-13"
$ws.Range("L10").Value = [double]"62.49850643161025"
$ws.Range("M10").Value = "[46.28679678539021, 78.7102160778303]"
$ws.Range("N10").Value = [double]"7.601714813176841e-10"
$ws.Range("O10").Value = [double]"7.601714813176841e-10"
$ws.Range("P10").Value = [double]"2.396289891927349"
$ws.Range("Q10").Value = "[2.132131951084965, 2.6604478327697336]"
$ws.Range("T10").Value = [double]"57.91732778796447"
$ws.Range("U10").Value = "[49.394596231851295, 66.44005934407765]"
$ws.Range("X10").Value = [double]"15.32936936936964"
$ws.Range("Y10").Value = [double]"14.28756756756782"
$ws.Range("Z10").Value = [double]"16.37117117117146"
$ws.Range("F11").Value = [double]"24.78000000000043"
$ws.Range("H11").Value = [double]"5.564548821723747e-12"
$ws.Range("I11").Value = [double]"5.564548821723747e-12"
$ws.Range("L11").Value = [double]"57.7675996773284"
$ws.Range("M11").Value = "[41.426776828680374, 74.10842252597642]"
$ws.Range("N11").Value = [double]"6.761792237242048e-09"
$ws.Range("O11").Value = [double]"6.761792237242048e-09"
$ws.Range("P11").Value = [double]"2.257921446724196"
$ws.Range("Q11").Value = "[1.9811845563178876, 2.5346583371305043]"
$ws.Range("T11").Value = [double]"53.17447688465406"
$ws.Range("U11").Value = "[44.764048700286125, 61.584905069021985]"
$ws.Range("V11").Value = [double]"2.220446049250313e-16"
$ws.Range("W11").Value = [double]"2.220446049250313e-16"
$ws.Range("X11").Value = [double]"15.87507507507535"
$ws.Range("Y11").Value = [double]"14.78366366366392"
$ws.Range("Z11").Value = [double]"16.96648648648678"
$ws.Range("F12").Value = [double]"24.78000000000043"
$ws.Range("H12").Value = [double]"1.023292561797007e-12"
$ws.Range("I12").Value = [double]"1.023292561797007e-12"
$ws.Range("L12").Value = [double]"59.3307576607934"
$ws.Range("M12").Value = "[44.59744555135808, 74.06406977022873]"
$ws.Range("N12").Value = [double]"2.379587638046132e-10"
$ws.Range("O12").Value = [double]"2.379587638046132e-10"
$ws.Range("P12").Value = [double]"1.993763505881811"
$ws.Range("Q12").Value = "[1.7296055650394253, 2.257921446724196]"
$ws.Range("T12").Value = [double]"54.11771106017719"
$ws.Range("U12").Value = "[45.889980564321526, 62.34544155603286]"
$ws.Range("X12").Value = [double]"16.91687687687718"
$ws.Range("Y12").Value = [double]"15.87507507507535"
$ws.Range("Z12").Value = [double]"17.958678678679"
$ws.Range("F13").Value = [double]"24.78000000000043"
$ws.Range("H13").Value = [double]"6.641354133307686e-13"
$ws.Range("I13").Value = [double]"6.641354133307686e-13"
$ws.Range("L13").Value = [double]"53.97796438485651"
$ws.Range("M13").Value = "[41.802974153636214, 66.15295461607681]"
$ws.Range("N13").Value = [double]"1.593725151849412e-11"
$ws.Range("O13").Value = [double]"1.593725151849412e-11"
$ws.Range("P13").Value = [double]"1.754763464167272"
$ws.Range("Q13").Value = "[1.4906055233248878, 2.0189214050096567]"
$ws.Range("T13").Value = [double]"55.1222170761482"
$ws.Range("U13").Value = "[47.717265198011155, 62.52716895428524]"
$ws.Range("X13").Value = [double]"17.85945945945977"
$ws.Range("Y13").Value = [double]"16.81765765765795"
$ws.Range("Z13").Value = [double]"18.90126126126159"
$ws.Range("F14").Value = [double]"24.78000000000043"
$ws.Range("H14").Value = [double]"3.639311074721263e-13"
$ws.Range("I14").Value = [double]"3.639311074721263e-13"
$ws.Range("L14").Value = [double]"57.41588495578753"
$ws.Range("M14").Value = "[45.78957786917246, 69.0421920424026]"
$ws.Range("N14").Value = [double]"6.148415110374117e-13"
$ws.Range("O14").Value = [double]"6.148415110374117e-13"
$ws.Range("P14").Value = [double]"1.628973968528041"
$ws.Range("Q14").Value = "[1.3899739268135027, 1.8679740102425795]"
$ws.Range("T14").Value = [double]"53.88636356747125"
$ws.Range("U14").Value = "[46.157369763832584, 61.615357371109916]"
$ws.Range("X14").Value = [double]"18.35555555555588"
$ws.Range("Y14").Value = [double]"17.41297297297328"
$ws.Range("Z14").Value = [double]"19.29813813813848"

Write-Output "applied"